$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (price & 1h volume %) per commit diff.
# Force each target cell to Text format before assignment so Excel
# stores the literal string instead of re-interpreting it as a number
# or percentage (matches the workbook's original inlineStr text cells).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.84"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "5.42%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "9.74%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.263"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.31%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07501"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "6.14%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.872"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.74%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.814"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "7.34%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.491"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "6.55%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9209"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.96%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1699"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "6.06%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07869"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.83%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08085"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "5.41%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03050"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.77%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09903"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "9.67%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001499"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-5.93%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04600"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.73%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006427"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3.24%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.462"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.62%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.229"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.13%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3303"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.86%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1342"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.32%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.498"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "12.38%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1619"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.41%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.51%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.83%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001397"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "19.65%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "16.04%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2,544.94%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04488"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.94%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006875"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.94%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1348"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "7.64%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002235"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "8.20%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01277"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "9.28%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006161"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.02%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7094"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-63.23%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "15.42%"
